$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / metadata text updates ---
$ws.Range("M6").Value = "Edward A. Caban"
$ws.Range("A8").Value = "Volume 30   Number  27"
$ws.Range("C9").Value = "Report Covering the Week  7/3/2023  Through  7/9/2023"

# --- Crime statistics table updates ---
$ws.Range("N15").Value = -66.666666666666
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = -36.363636363636
$ws.Range("I16").Value = 51
$ws.Range("J16").Value = 62
$ws.Range("K16").Value = -17.741935483871
$ws.Range("L16").Value = 4.081632653061
$ws.Range("M16").Value = 2
$ws.Range("N16").Value = -83.86075949367
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 300
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = -22.222222222222
$ws.Range("I17").Value = 59
$ws.Range("J17").Value = 55
$ws.Range("K17").Value = 7.272727272727
$ws.Range("L17").Value = 90.322580645161
$ws.Range("M17").Value = 55.263157894736
$ws.Range("N17").Value = 20.408163265306
$ws.Range("C18").Value = 5
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = -52.631578947368
$ws.Range("I18").Value = 62
$ws.Range("J18").Value = 79
$ws.Range("K18").Value = -21.518987341772
$ws.Range("L18").Value = 29.166666666666
$ws.Range("M18").Value = 5.084745762711
$ws.Range("N18").Value = -85.273159144893
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = -13.333333333333
$ws.Range("F19").Value = 58
$ws.Range("G19").Value = 73
$ws.Range("H19").Value = -20.547945205479
$ws.Range("I19").Value = 378
$ws.Range("J19").Value = 413
$ws.Range("K19").Value = -8.474576271186
$ws.Range("L19").Value = 19.620253164557
$ws.Range("M19").Value = 13.855421686747
$ws.Range("N19").Value = -62.977473065621
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 33.333333333333
$ws.Range("I20").Value = 49
$ws.Range("J20").Value = 29
$ws.Range("K20").Value = 68.965517241379
$ws.Range("L20").Value = 19.512195121951
$ws.Range("M20").Value = 172.222222222222
$ws.Range("N20").Value = -91.764705882352
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = 8.695652173913
$ws.Range("F21").Value = 89
$ws.Range("G21").Value = 118
$ws.Range("H21").Value = -24.57627118644
$ws.Range("I21").Value = 605
$ws.Range("J21").Value = 645
$ws.Range("K21").Value = -6.201550387596
$ws.Range("L21").Value = 23.469387755102
$ws.Range("M21").Value = 20.517928286852
$ws.Range("N21").Value = -74.979321753515
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -50
$ws.Range("I22").Value = 15
$ws.Range("J22").Value = 14
$ws.Range("K22").Value = 7.142857142857
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = -33.333333333333
$ws.Range("I23").Value = 24
$ws.Range("J23").Value = 13
$ws.Range("K23").Value = 84.615384615384
$ws.Range("L23").Value = 118.181818181818
$ws.Range("M23").Value = 41.176470588235
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = -16.666666666666
$ws.Range("F24").Value = 92
$ws.Range("G24").Value = 93
$ws.Range("H24").Value = -1.075268817204
$ws.Range("I24").Value = 546
$ws.Range("J24").Value = 712
$ws.Range("K24").Value = -23.314606741573
$ws.Range("L24").Value = -25.511596180081
$ws.Range("M24").Value = 4.798464491362
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 24
$ws.Range("G25").Value = 17
$ws.Range("H25").Value = 41.176470588235
$ws.Range("I25").Value = 122
$ws.Range("J25").Value = 121
$ws.Range("K25").Value = 0.826446280991
$ws.Range("L25").Value = 43.529411764705
$ws.Range("M25").Value = -14.084507042253
$ws.Range("J26").Value = 9
$ws.Range("K26").Value = 44.444444444444
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -42.857142857142
$ws.Range("I27").Value = 21
$ws.Range("J27").Value = 35
$ws.Range("K27").Value = -40
$ws.Range("L27").Value = -4.545454545454
$ws.Range("H30").Value = -100

# --- Cells changing from text placeholder to numeric value ---
$ws.Range("I28").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = 4
$ws.Range("L28").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E18").Value = 25
$ws.Range("I28").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C23").Value = 2
$ws.Range("I28").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("D26").Value = 1
$ws.Range("L28").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("E26").Value = -100

# --- Cells changing from numeric value to text placeholder ---
$ws.Range("C22").Value = "'0"
$ws.Range("C28").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D23").Value = "'0"
$ws.Range("C28").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E23").Value = "'***.*"
$ws.Range("C28").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("F30").Value = "'0"
$ws.Range("C28").Copy()
$ws.Range("F30").PasteSpecial(-4122)
